$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textCells = @("D5","D6","D7","D8","D9","D10","D11","D12","D13","D14","D15","D17","D18","D19","D20","D21","D22","D23","D25","D26","D27","D28","D29","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $textCells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range("D2").Value = '24.315.66'
$ws.Range("E2").Value = '  +8.83%  '
$ws.Range("D3").Value = '1.673.89'
$ws.Range("E3").Value = '  +4.18%  '
$ws.Range("E4").Value = '  -0.35%  '
$ws.Range("D5").Value = '307.19'
$ws.Range("E5").Value = '  +5.88%  '
$ws.Range("D6").Value = '0.9965'
$ws.Range("E6").Value = '  +0.51%  '
$ws.Range("D7").Value = '0.3706'
$ws.Range("E7").Value = '  -0.28%  '
$ws.Range("D8").Value = '0.3444'
$ws.Range("E8").Value = '  +1.77%  '
$ws.Range("D9").Value = '47.22'
$ws.Range("E9").Value = '  +9.90%  '
$ws.Range("D10").Value = '1.184'
$ws.Range("E10").Value = '  +2.62%  '
$ws.Range("D11").Value = '0.07252'
$ws.Range("E11").Value = '  +2.13%  '
$ws.Range("D12").Value = '0.9991'
$ws.Range("E12").Value = '  -0.07%  '
$ws.Range("D13").Value = '20.36'
$ws.Range("E13").Value = '  +1.35%  '
$ws.Range("D14").Value = '6.096'
$ws.Range("E14").Value = '  +2.25%  '
$ws.Range("D15").Value = '6.743'
$ws.Range("E15").Value = '  +1.39%  '
$ws.Range("D16").Value = '1.673.27'
$ws.Range("E16").Value = '  +4.35%  '
$ws.Range("D17").Value = '0.00001109'
$ws.Range("E17").Value = '  +2.11%  '
$ws.Range("D18").Value = '0.9969'
$ws.Range("E18").Value = '  +0.72%  '
$ws.Range("D19").Value = '0.06707'
$ws.Range("E19").Value = '  +0.14%  '
$ws.Range("D20").Value = '81.18'
$ws.Range("E20").Value = '  +3.75%  '
$ws.Range("D21").Value = '16.47'
$ws.Range("E21").Value = '  +1.37%  '
$ws.Range("D22").Value = '6.091'
$ws.Range("E22").Value = '  +0.06%  '
$ws.Range("D23").Value = '11.95'
$ws.Range("E23").Value = '  +1.06%  '
$ws.Range("D24").Value = '24.257.69'
$ws.Range("E24").Value = '  +8.23%  '
$ws.Range("D25").Value = '2.425'
$ws.Range("E25").Value = '  +1.05%  '
$ws.Range("D26").Value = '3.360'
$ws.Range("E26").Value = '  -9.85%  '
$ws.Range("D27").Value = '2.662'
$ws.Range("E27").Value = '  +4.69%  '
$ws.Range("D28").Value = '151.92'
$ws.Range("E28").Value = '  +1.00%  '
$ws.Range("D29").Value = '19.55'
$ws.Range("E29").Value = '  -1.03%  '
$ws.Range("D30").Value = '1.853.79'
$ws.Range("E30").Value = '  +3.80%  '
$ws.Range("D31").Value = '127.25'
$ws.Range("E31").Value = '  +4.50%  '
$ws.Range("D32").Value = '6.309'
$ws.Range("E32").Value = '  +5.14%  '
$ws.Range("D33").Value = '4.030'
$ws.Range("E33").Value = '  -6.37%  '
$ws.Range("D34").Value = '0.9711'
$ws.Range("E34").Value = '  +1.38%  '
$ws.Range("D35").Value = '1.754'
$ws.Range("E35").Value = '  +6.70%  '
$ws.Range("D36").Value = '0.08457'
$ws.Range("E36").Value = '  +2.31%  '
$ws.Range("D37").Value = '9.025'
$ws.Range("E37").Value = '  +3.41%  '
$ws.Range("D38").Value = '12.26'
$ws.Range("E38").Value = '  +3.67%  '
$ws.Range("B39").Value = 'InternetComputer(DFINITY)'
$ws.Range("C39").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D39").Value = '5.346'
$ws.Range("E39").Value = '  +0.45%  '
$ws.Range("B40").Value = 'Hedera'
$ws.Range("C40").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D40").Value = '0.06383'
$ws.Range("E40").Value = '  +2.50%  '
$ws.Range("D41").Value = '0.02341'
$ws.Range("E41").Value = '  +5.89%  '
$ws.Range("D42").Value = '1.258'
$ws.Range("E42").Value = '  +1.50%  '
$ws.Range("D43").Value = '0.2107'
$ws.Range("E43").Value = '  +4.18%  '
$ws.Range("D44").Value = '0.6161'
$ws.Range("E44").Value = '  +2.47%  '
$ws.Range("D45").Value = '0.9962'
$ws.Range("E45").Value = '  +0.71%  '
$ws.Range("B46").Value = 'PancakeSwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D46").Value = '3.783'
$ws.Range("E46").Value = '  +3.28%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = '13.07'
$ws.Range("E47").Value = '  -1.88%  '
$ws.Range("D48").Value = '0.5934'
$ws.Range("E48").Value = '  +2.89%  '
$ws.Range("D49").Value = '127.23'
$ws.Range("E49").Value = '  +1.73%  '
$ws.Range("D50").Value = '2.026'
$ws.Range("E50").Value = '  +2.58%  '
$ws.Range("D51").Value = '0.07209'
$ws.Range("E51").Value = '  +4.86%  '
